$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '36.394.93'
$ws.Range("E2").Value = '  -1.10%  '

# Row 3
$ws.Range("D3").Value = '2.036.30'
$ws.Range("E3").Value = '  -2.55%  '

# Row 4
$ws.Range("E4").Value = '  +0.34%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '244.87'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.58%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.660'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.29%  '

# Row 7
$ws.Range("E7").Value = '  +0.08%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '55.47'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.91%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '62.79'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +5.65%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.363'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.30%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0740'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -3.73%  '

# Row 12
$ws.Range("E12").Value = '  -3.18%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.893'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +1.38%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.342.76'
$ws.Range("E14").Value = '  -1.94%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '14.04'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -5.63%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '5.34'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -4.08%  '

# Row 17
$ws.Range("D17").Value = '2.040.64'
$ws.Range("E17").Value = '  -2.23%  '

# Row 18
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '36.325.35'
$ws.Range("E18").Value = '  -1.24%  '

# Row 19
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '17.29'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.25%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '71.22'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.60%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0849'
$ws.Range("E21").Value = '  -3.43%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '235.87'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.17%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.15'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -6.03%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.20%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.35'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.85%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.23'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +2.44%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.17'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -8.00%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '162.96'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -3.15%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '19.84'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -5.91%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.120'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -2.63%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.19'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.74%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.92'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -7.56%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.0594'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.72%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.36'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -7.76%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0881'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +4.45%  '

# Row 36
$ws.Range("E36").Value = '  +0.31%  '

# Row 37
$ws.Range("E37").Value = '  -0.77%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.17'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -9.34%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '5.03'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.59%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.21'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -5.78%  '

# Row 41
$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.87'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -2.20%  '

# Row 42
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0214'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -3.13%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.09'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -6.05%  '

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '92.66'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -4.01%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0900'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -5.94%  '

# Row 46
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.364.53'
$ws.Range("E46").Value = '  +2.22%  '

# Row 47
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '15.64'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -4.46%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.37'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +5.22%  '

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.94'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +2.12%  '

# Row 50
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.226.06'
$ws.Range("E50").Value = '  -2.06%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '45.45'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.52%  '
